$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (vitesse) values
$ws.Range("D2").Value = 500
$ws.Range("D3").Value = 46
$ws.Range("D4").Value = 850
$ws.Range("D5").Value = 750
$ws.Range("D6").Value = 54
$ws.Range("D7").Value = 18
$ws.Range("D8").Value = 34
$ws.Range("D9").Value = 54
$ws.Range("D10").Value = 83
$ws.Range("D12").Value = 23
$ws.Range("D13").Value = 18
$ws.Range("D14").Value = 54
$ws.Range("D15").Value = 103
$ws.Range("D16").Value = 21
$ws.Range("D17").Value = 16
$ws.Range("D18").Value = 158
$ws.Range("D19").Value = 12
$ws.Range("D20").Value = 3.6

# Update view: zoom and selection
$excel.ActiveWindow.Zoom = 78
$ws.Range("G12").Select()
